# Auswertung Teil 1 und 2 fertig
#
# 1) Fix the row-6 "Delta l" label on each Messung sheet: the generic
#    "$\Delta l_n$" placeholder becomes a per-sheet "$\Delta l_1$".."$\Delta l_4$".
# 2) Fix Messung1's row-6 shared formulas, which erroneously spilled into row 7
#    ("C6:K7" / "L4:L6") - they should match the other three sheets ("C6:K6" / "L6").
# 3) Add a new "Ergebnisse" sheet at the end that pulls the per-frequency
#    results together from Messung1..4 into one table.
# 4) Re-point the various sheet selections / active tab like the saved file did.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Messung1")
$ws2 = $wb.Worksheets.Item("Messung2")
$ws3 = $wb.Worksheets.Item("Messung3")
$ws4 = $wb.Worksheets.Item("Messung4")

# --- Messung1: repair the shared-formula ranges for row 6 (they incorrectly
#     extended to row 7 in the original file) and relabel A6. ---
$ws1.Range("C6:K6").Formula = "=C2-C3"
$ws1.Range("L6").Formula = "=MEDIAN(B6:K6)"
$ws1.Range("A6").Value = '$\Delta l_1$'

# --- Messung2 / Messung3 / Messung4: just relabel A6. ---
$ws2.Range("A6").Value = '$\Delta l_2$'
$ws3.Range("A6").Value = '$\Delta l_3$'
$ws4.Range("A6").Value = '$\Delta l_4$'

# --- Update each Messung sheet's selection (matches the saved workbook) and
#     leave them un-activated; the new results sheet ends up the active tab. ---
$ws1.Activate()
$ws1.Range("J22").Select()

$ws2.Activate()
$ws2.Range("A6").Select()

$ws3.Activate()
$ws3.Range("A6").Select()

$ws4.Activate()
$ws4.Range("A2").Select()

# --- Add the new "Ergebnisse" sheet after Messung4. ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $last)
$ws5.Name = "Ergebnisse"

# Header row
$ws5.Range("A1").Value = "Messung"
$ws5.Range("B1").Value = "Frequenz"
$ws5.Range("C1").Value = "T_v"
$ws5.Range("D1").Value = "T_n"
$ws5.Range("E1").Value = "n"
$ws5.Range("F1").Value = "Mittelwert Differenz"
$ws5.Range("G1").Value = "Lambda / cm"
$ws5.Range("H1").Value = "c"
$ws5.Range("I1").Value = "c_exp(20)"
$ws5.Range("J1").Value = "Mittelwert"
$ws5.Range("K1").Value = "Literatur"
$ws5.Range("L1").Value = "abs. Abweichung"
$ws5.Range("M1").Value = "rel. Abweichung"

# Row 2 -> Messung1
$ws5.Range("A2").Value = 1
$ws5.Range("B2").Formula = "=Messung1!P2"
$ws5.Range("C2").Formula = "=Messung1!M2"
$ws5.Range("D2").Formula = "=Messung1!N2"
$ws5.Range("E2").Formula = "=Messung1!`$O`$2"
$ws5.Range("F2").Formula = "=Messung1!`$L`$6"
$ws5.Range("G2").Formula = "=Messung1!M6"
$ws5.Range("H2").Formula = "=Messung1!N6"
$ws5.Range("I2").Formula = "=Messung1!O6"
$ws5.Range("J2").Formula = "=Messung1!P6"
$ws5.Range("K2").Formula = "=Messung1!Q6"
$ws5.Range("L2").Formula = "=Messung1!R6"
$ws5.Range("M2").Formula = "=Messung1!S6"

# Row 3 -> Messung2
$ws5.Range("A3").Value = 2
$ws5.Range("B3").Formula = "=Messung2!`$P`$2"
$ws5.Range("C3").Formula = "=Messung2!`$M`$2"
$ws5.Range("D3").Formula = "=Messung2!N2"
$ws5.Range("E3").Formula = "=Messung2!O2"
$ws5.Range("F3").Formula = "=Messung2!L6"
$ws5.Range("G3").Formula = "=Messung2!M6"
$ws5.Range("H3").Formula = "=Messung2!N6"
$ws5.Range("I3").Formula = "=Messung2!O6"
$ws5.Range("J3").Formula = "=Messung2!P6"
$ws5.Range("K3").Formula = "=Messung2!Q6"
$ws5.Range("L3").Formula = "=Messung2!R6"
$ws5.Range("M3").Formula = "=Messung2!S6"

# Row 4 -> Messung3
$ws5.Range("A4").Value = 3
$ws5.Range("B4").Formula = "=Messung3!`$P`$2"
$ws5.Range("C4").Formula = "=Messung3!M2"
$ws5.Range("D4").Formula = "=Messung3!N2"
$ws5.Range("E4").Formula = "=Messung3!O2"
$ws5.Range("F4").Formula = "=Messung3!L6"
$ws5.Range("G4").Formula = "=Messung3!M6"
$ws5.Range("H4").Formula = "=Messung3!N6"
$ws5.Range("I4").Formula = "=Messung3!O6"
$ws5.Range("J4").Formula = "=Messung3!P6"
$ws5.Range("K4").Formula = "=Messung3!Q6"
$ws5.Range("L4").Formula = "=Messung3!R6"
$ws5.Range("M4").Formula = "=Messung3!S6"

# Row 5 -> Messung4
$ws5.Range("A5").Value = 4
$ws5.Range("B5").Formula = "=Messung4!`$P`$2"
$ws5.Range("C5").Formula = "=Messung4!M2"
$ws5.Range("D5").Formula = "=Messung4!N2"
$ws5.Range("E5").Formula = "=Messung4!O2"
$ws5.Range("F5").Formula = "=Messung4!L6"
$ws5.Range("G5").Formula = "=Messung4!M6"
$ws5.Range("H5").Formula = "=Messung4!N6"
$ws5.Range("I5").Formula = "=Messung4!O6"
$ws5.Range("J5").Formula = "=Messung4!P6"
$ws5.Range("K5").Formula = "=Messung4!Q6"
$ws5.Range("L5").Formula = "=Messung4!R6"
$ws5.Range("M5").Formula = "=Messung4!S6"

# Column widths roughly matching the saved file's auto-fit columns.
$ws5.Columns.Item(1).AutoFit() | Out-Null
$ws5.Columns.Item(2).AutoFit() | Out-Null
$ws5.Columns.Item(3).AutoFit() | Out-Null
$ws5.Columns.Item(4).AutoFit() | Out-Null
$ws5.Columns.Item(5).AutoFit() | Out-Null
$ws5.Columns.Item(6).AutoFit() | Out-Null
$ws5.Columns.Item(7).AutoFit() | Out-Null
$ws5.Columns.Item(8).AutoFit() | Out-Null
$ws5.Columns.Item(9).AutoFit() | Out-Null
$ws5.Columns.Item(10).AutoFit() | Out-Null
$ws5.Columns.Item(11).AutoFit() | Out-Null
$ws5.Columns.Item(12).AutoFit() | Out-Null
$ws5.Columns.Item(13).AutoFit() | Out-Null

# Match the metric (2 cm) page margins used elsewhere in the workbook.
$ws5.PageSetup.TopMargin = 56.692913385826778
$ws5.PageSetup.BottomMargin = 56.692913385826778

# Ergebnisse becomes the active/selected sheet and cell, like the saved file.
$ws5.Activate()
$ws5.Range("D4").Select()
